# Update the dSF (column F) values to reflect the repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 0
    14 = 1
    19 = -7
    20 = -1
    25 = -14
    30 = -5
    31 = -1
    42 = 3
    51 = 3
    53 = -9
    54 = -10
    57 = -4
    62 = -1
    65 = -2
    66 = -1
    70 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
